# Update column F (dSF) values for rows 2-29, per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -5
    3  = -4
    4  = -3
    5  = 1
    6  = 1
    7  = -1
    8  = 10
    9  = 6
    10 = -2
    11 = 2
    12 = 2
    13 = -3
    14 = -4
    15 = 2
    16 = 2
    18 = 2
    19 = 3
    20 = -1
    21 = -3
    22 = -1
    23 = -2
    24 = 8
    25 = 1
    26 = 4
    27 = 4
    29 = -4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
